# Auto-generated edit script: apply scheduled runner price/profit updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 654222.2  # H28: 617980.25 -> 654222.2
$ws.Cells.Item(28, 9).Value = 1234842.8  # I28: 1852259.4 -> 1234842.8
$ws.Cells.Item(28, 10).Value = 1024  # J28: 840.75 -> 1024
$ws.Cells.Item(28, 11).Value = 1234842.8  # K28: 1852259.4 -> 1234842.8
$ws.Cells.Item(28, 12).Value = 1024  # L28: 840.75 -> 1024
$ws.Cells.Item(28, 13).Value = -1234357.8  # M28: -1851774.4 -> -1234357.8
$ws.Cells.Item(28, 14).Value = -1994  # N28: -1810.75 -> -1994

$ws.Cells.Item(62, 8).Value = 719215.0600000001  # H62: 1006200.1 -> 719215.0600000001
$ws.Cells.Item(62, 9).Value = 1431227.9  # I62: 2003477 -> 1431227.9
$ws.Cells.Item(62, 10).Value = 7202.2856  # J62: 8923.200000000001 -> 7202.2856
$ws.Cells.Item(62, 11).Value = 1431227.9  # K62: 2003477 -> 1431227.9
$ws.Cells.Item(62, 12).Value = 7202.2856  # L62: 8923.200000000001 -> 7202.2856
$ws.Cells.Item(62, 13).Value = -1430603.9  # M62: -2002853 -> -1430603.9
$ws.Cells.Item(62, 14).Value = -8450.285599999999  # N62: -10171.2 -> -8450.285599999999

$ws.Cells.Item(65, 8).Value = 719215.0600000001  # H65: 1006200.1 -> 719215.0600000001
$ws.Cells.Item(65, 9).Value = 1431227.9  # I65: 2003477 -> 1431227.9
$ws.Cells.Item(65, 10).Value = 7202.2856  # J65: 8923.200000000001 -> 7202.2856
$ws.Cells.Item(65, 11).Value = 7156139.5  # K65: 10017385 -> 7156139.5
$ws.Cells.Item(65, 12).Value = 36011.428  # L65: 44616 -> 36011.428
$ws.Cells.Item(65, 13).Value = -7153019.5  # M65: -10014265 -> -7153019.5
$ws.Cells.Item(65, 14).Value = -42251.428  # N65: -50856 -> -42251.428

$ws.Cells.Item(106, 8).Value = 22224424  # H106: 55557556 -> 22224424
$ws.Cells.Item(106, 9).Value = 22224424  # I106: 55557556 -> 22224424
$ws.Cells.Item(106, 11).Value = 22224424  # K106: 55557556 -> 22224424
$ws.Cells.Item(106, 13).Value = -22223793  # M106: -55556925 -> -22223793

$ws.Cells.Item(107, 8).Value = 463348.47  # H107: 585228.2 -> 463348.47
$ws.Cells.Item(107, 9).Value = 653932.9  # I107: 694795.9 -> 653932.9
$ws.Cells.Item(107, 10).Value = 500.57144  # J107: 867 -> 500.57144
$ws.Cells.Item(107, 11).Value = 653932.9  # K107: 694795.9 -> 653932.9
$ws.Cells.Item(107, 12).Value = 500.57144  # L107: 867 -> 500.57144
$ws.Cells.Item(107, 13).Value = -652012.9  # M107: -692875.9 -> -652012.9
$ws.Cells.Item(107, 14).Value = -4340.57144  # N107: -4707 -> -4340.57144

$ws.Cells.Item(113, 8).Value = 1933.1666  # H113: 2000 -> 1933.1666
$ws.Cells.Item(113, 9).Value = 1999.6666  # I113: 2000 -> 1999.6666
$ws.Cells.Item(113, 10).Value = 1866.6666  # J113: 2000 -> 1866.6666
$ws.Cells.Item(113, 11).Value = 1999.6666  # K113: 2000 -> 1999.6666
$ws.Cells.Item(113, 12).Value = 1866.6666  # L113: 2000 -> 1866.6666
$ws.Cells.Item(113, 13).Value = 1254.3334  # M113: 1254 -> 1254.3334
$ws.Cells.Item(113, 14).Value = -8374.6666  # N113: -8508 -> -8374.6666

$ws.Cells.Item(137, 8).Value = 1984.3334  # H137: 1886.7059 -> 1984.3334
$ws.Cells.Item(137, 9).Value = 1983.2142  # I137: 1969.1428 -> 1983.2142
$ws.Cells.Item(137, 10).Value = 2000  # J137: 1502 -> 2000
$ws.Cells.Item(137, 11).Value = 5949.642599999999  # K137: 5907.428400000001 -> 5949.642599999999
$ws.Cells.Item(137, 12).Value = 6000  # L137: 4506 -> 6000
$ws.Cells.Item(137, 13).Value = -3399.642599999999  # M137: -3357.428400000001 -> -3399.642599999999
$ws.Cells.Item(137, 14).Value = -11100  # N137: -9606 -> -11100

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2795.3215  # H32: 3474.6904 -> 2795.3215
$ws.Cells.Item(32, 9).Value = 1811.1489  # I32: 2248.147 -> 1811.1489
$ws.Cells.Item(32, 10).Value = 7934.8887  # J32: 8687.5 -> 7934.8887
$ws.Cells.Item(32, 11).Value = 1811.1489  # K32: 2248.147 -> 1811.1489
$ws.Cells.Item(32, 12).Value = 7934.8887  # L32: 8687.5 -> 7934.8887
$ws.Cells.Item(32, 13).Value = -1524.1489  # M32: -1961.147 -> -1524.1489
$ws.Cells.Item(32, 14).Value = -8508.8887  # N32: -9261.5 -> -8508.8887

$ws.Cells.Item(62, 8).Value = 0  # H62: 39900 -> 0
$ws.Cells.Item(62, 10).Value = 0  # J62: 39900 -> 0
$ws.Cells.Item(62, 12).Value = 0  # L62: 39900 -> 0
$ws.Cells.Item(62, 14).ClearContents()  # N62: -41148 -> (removed)

$ws.Cells.Item(65, 8).Value = 0  # H65: 39900 -> 0
$ws.Cells.Item(65, 10).Value = 0  # J65: 39900 -> 0
$ws.Cells.Item(65, 12).Value = 0  # L65: 119700 -> 0
$ws.Cells.Item(65, 14).ClearContents()  # N65: -125940 -> (removed)

$ws.Cells.Item(74, 8).Value = 8911.5625  # H74: 9463.467000000001 -> 8911.5625
$ws.Cells.Item(74, 9).Value = 1544.6923  # I74: 1684.909 -> 1544.6923
$ws.Cells.Item(74, 10).Value = 40834.668  # J74: 30854.5 -> 40834.668
$ws.Cells.Item(74, 11).Value = 1544.6923  # K74: 1684.909 -> 1544.6923
$ws.Cells.Item(74, 12).Value = 40834.668  # L74: 30854.5 -> 40834.668
$ws.Cells.Item(74, 13).Value = -670.6922999999999  # M74: -810.9090000000001 -> -670.6922999999999
$ws.Cells.Item(74, 14).Value = -42582.668  # N74: -32602.5 -> -42582.668

$ws.Cells.Item(77, 8).Value = 8911.5625  # H77: 9463.467000000001 -> 8911.5625
$ws.Cells.Item(77, 9).Value = 1544.6923  # I77: 1684.909 -> 1544.6923
$ws.Cells.Item(77, 10).Value = 40834.668  # J77: 30854.5 -> 40834.668
$ws.Cells.Item(77, 11).Value = 7723.461499999999  # K77: 8424.545 -> 7723.461499999999
$ws.Cells.Item(77, 12).Value = 204173.34  # L77: 154272.5 -> 204173.34
$ws.Cells.Item(77, 13).Value = -3355.461499999999  # M77: -4056.545 -> -3355.461499999999
$ws.Cells.Item(77, 14).Value = -212909.34  # N77: -163008.5 -> -212909.34

$ws.Cells.Item(132, 8).Value = 2365.0625  # H132: 2701.2903 -> 2365.0625
$ws.Cells.Item(132, 9).Value = 2024.4286  # I132: 2226.75 -> 2024.4286
$ws.Cells.Item(132, 10).Value = 4749.5  # J132: 4328.2856 -> 4749.5
$ws.Cells.Item(132, 11).Value = 6073.2858  # K132: 6680.25 -> 6073.2858
$ws.Cells.Item(132, 12).Value = 14248.5  # L132: 12984.8568 -> 14248.5
$ws.Cells.Item(132, 13).Value = -3543.2858  # M132: -4150.25 -> -3543.2858
$ws.Cells.Item(132, 14).Value = -19308.5  # N132: -18044.8568 -> -19308.5

$ws.Cells.Item(135, 8).Value = 31369.857  # H135: 50000 -> 31369.857
$ws.Cells.Item(135, 10).Value = 31369.857  # J135: 50000 -> 31369.857
$ws.Cells.Item(135, 12).Value = 31369.857  # L135: 50000 -> 31369.857
$ws.Cells.Item(135, 14).Value = -41509.857  # N135: -60140 -> -41509.857

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1112.1052  # H99: 1090.5264 -> 1112.1052
$ws.Cells.Item(99, 9).Value = 1118.3334  # I99: 1160 -> 1118.3334
$ws.Cells.Item(99, 10).Value = 1000  # J99: 500 -> 1000
$ws.Cells.Item(99, 11).Value = 1118.3334  # K99: 1160 -> 1118.3334
$ws.Cells.Item(99, 12).Value = 1000  # L99: 500 -> 1000
$ws.Cells.Item(99, 13).Value = 379.6666  # M99: 338 -> 379.6666
$ws.Cells.Item(99, 14).Value = -3996  # N99: -3496 -> -3996

$ws.Cells.Item(105, 8).Value = 2748  # H105: 2476.6316 -> 2748
$ws.Cells.Item(105, 9).Value = 2675.2258  # I105: 2686.9092 -> 2675.2258
$ws.Cells.Item(105, 10).Value = 3500  # J105: 2187.5 -> 3500
$ws.Cells.Item(105, 11).Value = 2675.2258  # K105: 2686.9092 -> 2675.2258
$ws.Cells.Item(105, 12).Value = 3500  # L105: 2187.5 -> 3500
$ws.Cells.Item(105, 13).Value = -928.2258000000002  # M105: -939.9092000000001 -> -928.2258000000002
$ws.Cells.Item(105, 14).Value = -6994  # N105: -5681.5 -> -6994

$ws.Cells.Item(123, 8).Value = 28975  # H123: 21240 -> 28975
$ws.Cells.Item(123, 10).Value = 29966.666  # J123: 20560 -> 29966.666
$ws.Cells.Item(123, 12).Value = 29966.666  # L123: 20560 -> 29966.666
$ws.Cells.Item(123, 14).Value = -39766.666  # N123: -30360 -> -39766.666

$ws.Cells.Item(134, 8).Value = 2325.5  # H134: 2693.9583 -> 2325.5
$ws.Cells.Item(134, 9).Value = 1685.0385  # I134: 1935.05 -> 1685.0385
$ws.Cells.Item(134, 11).Value = 5055.1155  # K134: 5805.15 -> 5055.1155
$ws.Cells.Item(134, 13).Value = -2520.1155  # M134: -3270.15 -> -2520.1155

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 5209226  # H99: 20834334 -> 5209226
$ws.Cells.Item(99, 9).Value = 6945311  # I99: 20834334 -> 6945311
$ws.Cells.Item(99, 10).Value = 971.3333  # J99: 0 -> 971.3333
$ws.Cells.Item(99, 11).Value = 6945311  # K99: 20834334 -> 6945311
$ws.Cells.Item(99, 12).Value = 971.3333  # L99: 0 -> 971.3333
$ws.Cells.Item(99, 13).Value = -6943813  # M99: -20832836 -> -6943813
$ws.Cells.Item(99, 14).Value = -3967.3333  # N99: None -> -3967.3333

$ws.Cells.Item(105, 8).Value = 539.8333  # H105: 496.66666 -> 539.8333
$ws.Cells.Item(105, 9).Value = 539.8333  # I105: 556 -> 539.8333
$ws.Cells.Item(105, 10).Value = 0  # J105: 200 -> 0
$ws.Cells.Item(105, 11).Value = 539.8333  # K105: 556 -> 539.8333
$ws.Cells.Item(105, 12).Value = 0  # L105: 200 -> 0
$ws.Cells.Item(105, 13).Value = 1207.1667  # M105: 1191 -> 1207.1667
$ws.Cells.Item(105, 14).ClearContents()  # N105: -3694 -> (removed)

$ws.Cells.Item(126, 8).Value = 5209226  # H126: 20834334 -> 5209226
$ws.Cells.Item(126, 9).Value = 6945311  # I126: 20834334 -> 6945311
$ws.Cells.Item(126, 10).Value = 971.3333  # J126: 0 -> 971.3333
$ws.Cells.Item(126, 11).Value = 20835933  # K126: 62503002 -> 20835933
$ws.Cells.Item(126, 12).Value = 2913.9999  # L126: 0 -> 2913.9999
$ws.Cells.Item(126, 13).Value = -20833463  # M126: -62500532 -> -20833463
$ws.Cells.Item(126, 14).Value = -7853.9999  # N126: None -> -7853.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 15443.75  # H68: 15631.25 -> 15443.75
$ws.Cells.Item(68, 9).Value = 60001  # I68: 24540.4 -> 60001
$ws.Cells.Item(68, 10).Value = 591.3333  # J68: 782.6667 -> 591.3333
$ws.Cells.Item(68, 11).Value = 180003  # K68: 73621.20000000001 -> 180003
$ws.Cells.Item(68, 12).Value = 1773.9999  # L68: 2348.0001 -> 1773.9999
$ws.Cells.Item(68, 13).Value = -179192  # M68: -72810.20000000001 -> -179192
$ws.Cells.Item(68, 14).Value = -3395.9999  # N68: -3970.0001 -> -3395.9999

$ws.Cells.Item(71, 8).Value = 15443.75  # H71: 15631.25 -> 15443.75
$ws.Cells.Item(71, 9).Value = 60001  # I71: 24540.4 -> 60001
$ws.Cells.Item(71, 10).Value = 591.3333  # J71: 782.6667 -> 591.3333
$ws.Cells.Item(71, 11).Value = 540009  # K71: 220863.6 -> 540009
$ws.Cells.Item(71, 12).Value = 5321.9997  # L71: 7044.0003 -> 5321.9997
$ws.Cells.Item(71, 13).Value = -535953  # M71: -216807.6 -> -535953
$ws.Cells.Item(71, 14).Value = -13433.9997  # N71: -15156.0003 -> -13433.9997

$ws.Cells.Item(131, 8).Value = 2515.7317  # H131: 2574.8354 -> 2515.7317
$ws.Cells.Item(131, 10).Value = 2600.1265  # J131: 2664.8948 -> 2600.1265
$ws.Cells.Item(131, 12).Value = 7800.379499999999  # L131: 7994.6844 -> 7800.379499999999
$ws.Cells.Item(131, 14).Value = -17880.3795  # N131: -18074.6844 -> -17880.3795

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 696487.75  # H122: 928059 -> 696487.75
$ws.Cells.Item(122, 9).Value = 795786.7  # I122: 928059 -> 795786.7
$ws.Cells.Item(122, 10).Value = 1395  # J122: 0 -> 1395
$ws.Cells.Item(122, 11).Value = 2387360.1  # K122: 2784177 -> 2387360.1
$ws.Cells.Item(122, 12).Value = 4185  # L122: 0 -> 4185
$ws.Cells.Item(122, 13).Value = -2384910.1  # M122: -2781727 -> -2384910.1
$ws.Cells.Item(122, 14).Value = -9085  # N122: None -> -9085

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(62, 8).Value = 0  # H62: 50000 -> 0
$ws.Cells.Item(62, 10).Value = 0  # J62: 50000 -> 0
$ws.Cells.Item(62, 12).Value = 0  # L62: 50000 -> 0
$ws.Cells.Item(62, 14).ClearContents()  # N62: -51248 -> (removed)

$ws.Cells.Item(65, 8).Value = 0  # H65: 50000 -> 0
$ws.Cells.Item(65, 10).Value = 0  # J65: 50000 -> 0
$ws.Cells.Item(65, 12).Value = 0  # L65: 150000 -> 0
$ws.Cells.Item(65, 14).ClearContents()  # N65: -156240 -> (removed)
